$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns T and U
$ws.Range("T1").Value = "P_discharge"
$ws.Range("U1").Value = "P_discharge_variance"

# Row 2 updates
$ws.Range("B2").Value = 61.9
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 9
$ws.Range("N2").Value = 51.6
$ws.Range("O2").Value = 1
$ws.Range("Q2").Value = 50
$ws.Range("R2").Value = 3
$ws.Range("T2").Value = 20
$ws.Range("U2").Value = 1

# Row 3 updates
$ws.Range("B3").Value = 47.5
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 7.2
$ws.Range("N3").Value = 95
$ws.Range("O3").Value = 1
$ws.Range("Q3").Value = 50
$ws.Range("R3").Value = 3
$ws.Range("T3").Value = 30
$ws.Range("U3").Value = 1

# Selection / view changes
$ws.Range("R9").Select()
